$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the data values in row 5 to 2 decimal places ("custom accuracy")
for ($col = 2; $col -le 34; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $v = $cell.Value()
    $cell.Value = [Math]::Round($v, 2)
}

# Narrow a few columns (K, Q, AC) to match the width used by other
# "7-wide" columns elsewhere in the sheet (e.g. column J / col 10)
$refWidth = $ws.Columns(10).ColumnWidth
$ws.Columns(11).ColumnWidth = $refWidth
$ws.Columns(17).ColumnWidth = $refWidth
$ws.Columns(28).ColumnWidth = $refWidth

# Remove the last data row (row 6) - trims the dataset down
$ws.Rows(6).Delete()
